$wb = $excel.ActiveWorkbook
Write-Output $excel.Windows.Count
$win = $excel.Windows.Item(1)
Write-Output $win.Left
